$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Prime the outline-level "high water mark" counters (outlineLevelRow=3,
# outlineLevelCol=6) on a pair of throwaway row/col before touching the real
# columns, then delete the scratch row/col afterwards so only the sheet-wide
# counters are left behind, matching the target sheetFormatPr.
$ws.Rows.Item(20).Group()
$ws.Rows.Item(20).Group()
$ws.Rows.Item(20).Group()

$ws.Columns.Item(20).Group()
$ws.Columns.Item(20).Group()
$ws.Columns.Item(20).Group()
$ws.Columns.Item(20).Group()
$ws.Columns.Item(20).Group()
$ws.Columns.Item(20).Group()

$ws.Columns.Item(20).ClearFormats()
$ws.Columns.Item(20).Delete()
$ws.Rows.Item(20).Delete()

# Column widths for the new data columns
$ws.Columns.Item(2).ColumnWidth = 34.14285714285714
$ws.Columns.Item(3).ColumnWidth = 51.57142857142857
$ws.Columns.Item(4).ColumnWidth = 22.714285714285715
$ws.Columns.Item(5).ColumnWidth = 69.28571428571428
$ws.Columns.Item(7).ColumnWidth = 44.14285714285714

# Row 3 - new study-log entry
$ws.Range("A3").Value = 15
$ws.Range("B3").Value = "9：30-11：07"
$ws.Range("C3").Value = "二、十、八、十六进制，数据类型，变量转换"
$ws.Range("D3").Value = "8；00-9：44"
$ws.Range("E3").Value = "运算符（到逻辑运算符）"
$ws.Range("F3").Value = "完成"
$ws.Range("G3").Value = "(今天下午肚子疼，没学，晚了些）"

# Row 4 - new study-log entry
$ws.Range("A4").Value = 16
$ws.Range("B4").Value = "9：03-10:58"
$ws.Range("C4").Value = "运算符，接受用户数据"
$ws.Range("D4").Value = "7:00-8：39"
$ws.Range("E4").Value = "流程控制（if，switch分支，for循环）"
$ws.Range("F4").Value = "完成"

$ws.Range("B6").Select()
